{"js": "// Remove the stray \"stop\" text that was typed between the tab stops and\n// the \"(podpis klienta)\" label in the signature block, restoring the\n// original blank line of tabs before that label.\nconst body = context.document.body;\nconst results = body.search(\"stop\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Remove the stray \"stop\" text that was typed between the tab stops and\n# the \"(podpis klienta)\" label in the signature block, restoring the\n# original blank line of tabs before that label.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\"stop\", $true, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n"}
